$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 used to hold "NA"; that value moves down to the newly appended row 3,
# so C2 becomes an empty (but still present) text cell. Writing a leading
# apostrophe keeps it text-typed instead of Excel coercing "" to a blank
# cell, and resetting the style back to "Normal" drops the quote-prefix
# formatting that the apostrophe would otherwise leave behind.
$ws.Range("C2").Value = "'"
$ws.Range("C2").Style = "Normal"

# New row 3: same term/count as row 2, with this run's date and the "NA"
# page-number value that used to live in C2.
$ws.Range("A3").Value = "'2025-03-03"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C3").Value = "NA"
$ws.Range("D3").Value = 1
